$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell $ws.Range("D2") "37.224.28"
Set-TextCell $ws.Range("E2") "  +1.70%  "

Set-TextCell $ws.Range("D3") "2.057.81"
Set-TextCell $ws.Range("E3") "  +1.02%  "

Set-TextCell $ws.Range("E4") "  +0.01%  "

Set-TextCell $ws.Range("D5") "232.83"
Set-TextCell $ws.Range("E5") "  +0.47%  "

Set-TextCell $ws.Range("D6") "0.620"
Set-TextCell $ws.Range("E6") "  +3.03%  "

Set-TextCell $ws.Range("D8") "57.27"
Set-TextCell $ws.Range("E8") "  +3.71%  "

Set-TextCell $ws.Range("E9") "  +3.44%  "

Set-TextCell $ws.Range("D10") "57.83"
Set-TextCell $ws.Range("E10") "  +1.64%  "

Set-TextCell $ws.Range("D11") "0.0759"
Set-TextCell $ws.Range("E11") "  +0.73%  "

Set-TextCell $ws.Range("E12") "  +1.07%  "

Set-TextCell $ws.Range("D13") "14.78"
Set-TextCell $ws.Range("E13") "  +3.23%  "

Set-TextCell $ws.Range("D14") "2.361.21"
Set-TextCell $ws.Range("E14") "  +1.08%  "

Set-TextCell $ws.Range("E15") "  +4.29%  "

Set-TextCell $ws.Range("D16") "0.783"
Set-TextCell $ws.Range("E16") "  +3.01%  "

Set-TextCell $ws.Range("D17") "5.17"
Set-TextCell $ws.Range("E17") "  -0.03%  "

Set-TextCell $ws.Range("D18") "2.058.55"
Set-TextCell $ws.Range("E18") "  +1.03%  "

Set-TextCell $ws.Range("D19") "37.186.37"
Set-TextCell $ws.Range("E19") "  +1.26%  "

Set-TextCell $ws.Range("D20") "6.36"
Set-TextCell $ws.Range("E20") "  +8.97%  "

Set-TextCell $ws.Range("D21") "69.26"
Set-TextCell $ws.Range("E21") "  +2.43%  "

Set-TextCell $ws.Range("D22") "0.0₃0810"
Set-TextCell $ws.Range("E22") "  +1.61%  "

Set-TextCell $ws.Range("D23") "225.71"
Set-TextCell $ws.Range("E23") "  +2.27%  "

Set-TextCell $ws.Range("E24") "  -0.01%  "

Set-TextCell $ws.Range("D25") "2.41"
Set-TextCell $ws.Range("E25") "  +0.30%  "

Set-TextCell $ws.Range("D26") "2.40"
Set-TextCell $ws.Range("E26") "  +1.28%  "

Set-TextCell $ws.Range("D27") "165.86"
Set-TextCell $ws.Range("E27") "  +1.86%  "

Set-TextCell $ws.Range("D28") "1.45"
Set-TextCell $ws.Range("E28") "  +7.71%  "

Set-TextCell $ws.Range("D29") "8.80"
Set-TextCell $ws.Range("E29") "  +0.84%  "

Set-TextCell $ws.Range("B30") "EthereumClassic"
Set-TextCell $ws.Range("C30") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D30") "19.03"
Set-TextCell $ws.Range("E30") "  +0.45%  "

Set-TextCell $ws.Range("B31") "Kaspa"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D31") "0.127"
Set-TextCell $ws.Range("E31") "  +0.28%  "

Set-TextCell $ws.Range("D32") "0.118"
Set-TextCell $ws.Range("E32") "  +0.78%  "

Set-TextCell $ws.Range("D33") "4.45"
Set-TextCell $ws.Range("E33") "  +2.13%  "

Set-TextCell $ws.Range("B34") "Hedera"
Set-TextCell $ws.Range("C34") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Range("D34") "0.0618"
Set-TextCell $ws.Range("E34") "  +2.15%  "

Set-TextCell $ws.Range("B35") "InternetComputer(DFINITY)"
Set-TextCell $ws.Range("C35") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D35") "4.61"
Set-TextCell $ws.Range("E35") "  +7.87%  "

Set-TextCell $ws.Range("E36") "  +0.48%  "

Set-TextCell $ws.Range("E37") "  +0.10%  "

Set-TextCell $ws.Range("E38") "  +1.62%  "

Set-TextCell $ws.Range("E39") "  -0.79%  "

Set-TextCell $ws.Range("D40") "5.69"
Set-TextCell $ws.Range("E40") "  -1.66%  "

Set-TextCell $ws.Range("E41") "  +0.31%  "

Set-TextCell $ws.Range("B42") "FTXToken"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell $ws.Range("D42") "4.42"
Set-TextCell $ws.Range("E42") "  -1.93%  "

Set-TextCell $ws.Range("B43") "Maker"
Set-TextCell $ws.Range("C43") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D43") "1.473.04"
Set-TextCell $ws.Range("E43") "  -0.17%  "

Set-TextCell $ws.Range("D44") "96.49"
Set-TextCell $ws.Range("E44") "  +3.04%  "

Set-TextCell $ws.Range("E45") "  +5.22%  "

Set-TextCell $ws.Range("D46") "0.0931"
Set-TextCell $ws.Range("E46") "  -1.33%  "

Set-TextCell $ws.Range("E47") "  +3.52%  "

Set-TextCell $ws.Range("E48") "  +1.79%  "

Set-TextCell $ws.Range("E49") "  -3.12%  "

Set-TextCell $ws.Range("E50") "  +3.15%  "

Set-TextCell $ws.Range("D51") "2.95"
Set-TextCell $ws.Range("E51") "  +1.82%  "
